$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 326; this shifts the existing rows 326..397
# down to 327..398, matching the target dimension change (A1:R397 -> A1:R398).
$ws.Rows.Item(326).Insert()

# Populate the newly inserted row 326 with the new record.
$ws.Range("A326").Value = 3
$ws.Range("B326").Value = "Femacal de La Calera"
$ws.Range("C326").Value = "Coquimbo"
$ws.Range("D326").Value = 44798
$ws.Range("D326").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E326").Value = 5
$ws.Range("F326").Value = 100112009
$ws.Range("G326").Value = "Acelga"
$ws.Range("H326").Value = "Sin especificar"
$ws.Range("I326").Value = "Primera"
$ws.Range("J326").Value = 280
$ws.Range("K326").Value = 3200
$ws.Range("L326").Value = 3300
$ws.Range("M326").Value = 3257
$ws.Range("N326").Value = "$/docena de atados (6 kilos)"
$ws.Range("O326").Value = "Provincia de Quillota"
$ws.Range("P326").Value = 543
$ws.Range("Q326").Value = 6
$ws.Range("R326").Value = "Hortaliza"
